# Adds 16 new match rows (rows 129-144) to the Liga Brasil 2025 sheet,
# corresponding to matchdays played 2025-07-19 through 2025-07-24.
# Columns: Fecha|Local|Visita|Goles Local|Goles Visita|Fixture ID|
#          Corners Local|Corners Visita|Amarillas Local|Amarillas Visita|
#          Rojas Local|Rojas Visita|Goles 1T Local|Goles 1T Visita|
#          Goles 2T Local|Goles 2T Visita|Posesion Local (%)|Posesion Visita (%)|Resultado

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 129

$matches = @(
    "2025-07-19|Fortaleza EC|Bahia|1|1|1351191|9|3|5|3|0|0|0|0|1|1|35%|65%|E",
    "2025-07-19|Vasco DA Gama|Gremio|1|1|1351184|12|2|2|3|0|0|0|0|1|1|64%|36%|E",
    "2025-07-19|Mirassol|Santos|3|0|1351187|9|2|1|1|1|0|0|0|3|0|43%|57%|L",
    "2025-07-20|Sao Paulo|Corinthians|2|0|1351186|6|5|2|2|0|0|0|0|2|0|48%|52%|L",
    "2025-07-20|Internacional|Ceara|1|0|1351189|7|5|1|2|0|0|0|0|1|0|46%|54%|L",
    "2025-07-20|Cruzeiro|Juventude|4|0|1351188|8|3|1|2|0|0|0|0|4|0|48%|52%|L",
    "2025-07-20|Vitoria|RB Bragantino|1|0|1351190|2|5|2|2|0|0|0|0|1|0|28%|72%|L",
    "2025-07-20|Palmeiras|Atletico-MG|3|2|1351185|9|0|2|3|0|0|0|0|3|2|58%|42%|L",
    "2025-07-20|Sport Recife|Botafogo|0|1|1351192|6|4|1|2|0|0|0|0|0|1|41%|59%|V",
    "2025-07-20|Flamengo|Fluminense|1|0|1351183|8|3|3|1|0|0|0|0|1|0|66%|34%|L",
    "2025-07-23|Fluminense|Palmeiras|1|2|1351193|5|3|1|2|0|1|0|0|1|2|54%|46%|V",
    "2025-07-23|Ceara|Mirassol|0|2|1351201|14|6|2|3|0|0|0|0|0|2|56%|44%|V",
    "2025-07-23|Corinthians|Cruzeiro|0|0|1351195|4|2|4|2|0|0|0|0|0|0|47%|53%|E",
    "2025-07-24|Santos|Internacional|1|2|1351196|8|2|4|2|0|0|0|0|1|2|65%|35%|V",
    "2025-07-24|Vitoria|Sport Recife|2|2|1351200|6|10|4|3|0|0|0|0|2|2|48%|52%|E",
    "2025-07-24|RB Bragantino|Flamengo|1|2|1351197|2|3|1|5|0|0|0|0|1|2|37%|63%|V"
)

$endRow = $startRow + $matches.Length - 1

# Pre-format the date column and the two possession-percentage columns as
# Text so Excel stores the values exactly as literal strings (e.g. "35%")
# instead of auto-converting them into percentage numbers / date serials.
$dateRange = $ws.Range("A" + $startRow + ":A" + $endRow)
$dateRange.NumberFormat = "@"
$pctRange = $ws.Range("Q" + $startRow + ":R" + $endRow)
$pctRange.NumberFormat = "@"

for ($i = 0; $i -lt $matches.Length; $i++) {
    $r = $startRow + $i
    $fields = $matches[$i].Split("|")

    $ws.Cells.Item($r, 1).Value = $fields[0]          # Fecha
    $ws.Cells.Item($r, 2).Value = $fields[1]          # Local
    $ws.Cells.Item($r, 3).Value = $fields[2]          # Visita
    $ws.Cells.Item($r, 4).Value = [int]$fields[3]     # Goles Local
    $ws.Cells.Item($r, 5).Value = [int]$fields[4]     # Goles Visita
    $ws.Cells.Item($r, 6).Value = [int]$fields[5]     # Fixture ID
    $ws.Cells.Item($r, 7).Value = [int]$fields[6]     # Corners Local
    $ws.Cells.Item($r, 8).Value = [int]$fields[7]     # Corners Visita
    $ws.Cells.Item($r, 9).Value = [int]$fields[8]     # Amarillas Local
    $ws.Cells.Item($r, 10).Value = [int]$fields[9]    # Amarillas Visita
    $ws.Cells.Item($r, 11).Value = [int]$fields[10]   # Rojas Local
    $ws.Cells.Item($r, 12).Value = [int]$fields[11]   # Rojas Visita
    $ws.Cells.Item($r, 13).Value = [int]$fields[12]   # Goles 1T Local
    $ws.Cells.Item($r, 14).Value = [int]$fields[13]   # Goles 1T Visita
    $ws.Cells.Item($r, 15).Value = [int]$fields[14]   # Goles 2T Local
    $ws.Cells.Item($r, 16).Value = [int]$fields[15]   # Goles 2T Visita
    $ws.Cells.Item($r, 17).Value = $fields[16]        # Posesion Local (%)
    $ws.Cells.Item($r, 18).Value = $fields[17]        # Posesion Visita (%)
    $ws.Cells.Item($r, 19).Value = $fields[18]        # Resultado
}

# Restore the default (Normal) style on the text-formatted ranges so the
# cells don't carry a lingering explicit number format, matching the rest
# of the sheet's plain/unstyled data cells.
$dateRange.Style = "Normal"
$pctRange.Style = "Normal"
